# Updated RAD Test Cases.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Existing Liability w/Notice Number
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Fri Sep 08 18:11:18 EDT 2023"
$ws.Range("C2").Value = "Y"
$ws.Range("D2").Value = "Existing Liability w/Notice Number"
$ws.Range("E2").Value = "Personal Income Tax"

# Row 3: Quarterly Estimated Tax
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Sep 08 18:11:31 EDT 2023"
$ws.Range("C3").Value = "Y"
$ws.Range("D3").Value = "Quarterly Estimated Tax"
$ws.Range("E3").Value = "Personal Income Tax"

# Row 4: Extension Payments
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Sep 08 18:11:45 EDT 2023"
$ws.Range("C4").Value = "Y"
$ws.Range("D4").Value = "Extension Payments"
$ws.Range("E4").Value = "Personal Income Tax"

# Row 5: New Tax Return Amount Due
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Sep 08 18:11:58 EDT 2023"
$ws.Range("C5").Value = "Y"
$ws.Range("D5").Value = "New Tax Return Amount Due"
$ws.Range("E5").Value = "Personal Income Tax"

# Update the selected cell to match the saved sheet view state
$ws.Range("D4").Select()
